$wb = $excel.ActiveWorkbook

# --- Rushing sheet: remove A.Brown's row ---
$ws1 = $wb.Worksheets.Item("Rushing")
$ws1.Rows.Item(9).Delete()

# --- Receiving sheet: remove A.Brown's row, add J.Jones at the end ---
$ws2 = $wb.Worksheets.Item("Receiving")
$ws2.Rows.Item(6).Delete()

$ws2.Cells.Item(15, 1).Value = 14
$ws2.Cells.Item(15, 2).Value = "J.Jones"
$ws2.Cells.Item(15, 3).Value = 42
$ws2.Cells.Item(15, 4).Value = 35
$ws2.Cells.Item(15, 5).Value = 9
$ws2.Cells.Item(15, 6).Value = 8
$ws2.Cells.Item(15, 7).Value = 5
$ws2.Cells.Item(15, 8).Value = 4

# Give the new last row's A cell the same look as the rest of the column
# (bold, centered, thin left/right border) but without a bottom border,
# matching how the bottom edge of the table moved down by one row.
$ws2.Cells.Item(14, 1).Copy()
$ws2.Cells.Item(15, 1).PasteSpecial(-4122)
$ws2.Cells.Item(15, 1).Value = 14
$ws2.Cells.Item(15, 1).Borders.Item(9).LineStyle = -4142
$ws2.Cells.Item(15, 1).Borders.Item(8).LineStyle = -4142
